$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters whose values change (G, L, M, N, O, A remain untouched)
$cols = @("B","C","D","E","F","H","I","J","K")

# New loading-percent values per row (row 2 .. row 25), in column order above
$data = @{
    2 = @(12.82367325436917, 10.80812076374232, 5.76211399874259, 12.26887857873896, 49.62168194613091, 7.344005520526261, 33.58131821716631, 10.53500721825178, 14.33777012433201)
    3 = @(12.70464726548083, 10.67941146970937, 5.784470725151469, 12.25277832548053, 49.3387827644138, 7.344005520526261, 33.44895099559951, 10.53941767323601, 14.26200877368644)
    4 = @(12.63556066275467, 10.60304012658475, 5.799911495454372, 12.2459384384422, 49.17254849473664, 7.344005520526261, 33.37163555900018, 10.54381052119563, 14.21964016793561)
    5 = @(12.60844648370843, 10.57262339028261, 5.806634368111398, 12.24391928002736, 49.10672361304903, 7.344005520526261, 33.34113651586447, 10.54602432384748, 14.20343390797324)
    6 = @(12.6040079184835, 10.5676163612037, 5.807776695825363, 12.2436304490014, 49.09591030950764, 7.344005520526261, 33.33613336489238, 10.54641751159059, 14.20080728168058)
    7 = @(12.63519074017636, 10.60262701260498, 5.800000419265411, 12.24590809460755, 49.17165294624712, 7.344005520526261, 33.37122014326444, 10.54383866190991, 14.21941729586253)
    8 = @(12.78182507176832, 10.76321317381672, 5.769466781920062, 12.26269661963071, 49.52261298582634, 7.344005520526261, 33.5348640517654, 10.53617830776234, 14.31079653714764)
    9 = @(13.09946322363903, 11.09749836741337, 5.723197490012542, 12.31965864233347, 50.26823920801485, 7.344005520526261, 33.8865981869458, 10.534520002061, 14.52213189331218)
    10 = @(13.34888342351625, 11.3525554958027, 5.697508689806755, 12.37595531932182, 50.84839111424846, 7.344005520526261, 34.16297236067171, 10.54143688829296, 14.69584022476559)
    11 = @(13.46532206065583, 11.47012957071416, 5.687626214346906, 12.40464562244898, 51.11872631340367, 7.344005520526261, 34.29240080880835, 10.54634492011859, 14.77859378814999)
    12 = @(13.50979812393752, 11.51483145700966, 5.684143231824011, 12.41594701778058, 51.22196350732769, 7.344005520526261, 34.34192600303662, 10.54845605109824, 14.81044227488538)
    13 = @(13.50020307610806, 11.50519682465429, 5.684881823784922, 12.41349372621859, 51.19969187823919, 7.344005520526261, 34.3312373464147, 10.54799016160064, 14.80356079658742)
    14 = @(13.4689737087073, 11.4738038785141, 5.687334471155041, 12.40556667337917, 51.127202649197, 7.344005520526261, 34.29646509946354, 10.54651354707115, 14.78120387912789)
    15 = @(13.44989340199574, 11.45459684143361, 5.688870553743785, 12.40076784677548, 51.08291213005906, 7.344005520526261, 34.27523237399387, 10.54564194523873, 14.76757546846945)
    16 = @(13.34132996057756, 11.34489917373827, 5.698190821969066, 12.37414183263199, 50.83084851617996, 7.344005520526261, 34.15458682578272, 10.54115152320783, 14.69050506165092)
    17 = @(13.27546039875709, 11.27796854807897, 5.704370418745464, 12.35859259368253, 50.67782082162494, 7.344005520526261, 34.08151008609278, 10.538847480401, 14.64416327253779)
    18 = @(13.23785680620009, 11.23962004424559, 5.708094520389976, 12.34993944710568, 50.5904112167123, 7.344005520526261, 34.03982817041989, 10.53768811043357, 14.61786162931312)
    19 = @(13.22517480871217, 11.22666268415799, 5.709384593682311, 12.34705968346387, 50.56092191466207, 7.344005520526261, 34.02577605566789, 10.53732407341582, 14.60901771337336)
    20 = @(13.2824433761639, 11.28507840244622, 5.703695020556138, 12.36021782663762, 50.69404837011796, 7.344005520526261, 34.08925315400219, 10.53907558927601, 14.64906008822743)
    21 = @(13.47813648647527, 11.48302022954024, 5.686607032723112, 12.40788323164412, 51.14847141603823, 7.344005520526261, 34.30666477420658, 10.54694041657399, 14.78775697173695)
    22 = @(13.60824641611644, 11.61341104720637, 5.676950339485074, 12.44157930171674, 51.45049017528974, 7.344005520526261, 34.45173907222832, 10.55355223919295, 14.88137248592563)
    23 = @(13.53861661795755, 11.54373948348968, 5.681966052215012, 12.42336443565161, 51.28885610658492, 7.344005520526261, 34.37404380828264, 10.54988900472107, 14.83114504497181)
    24 = @(13.279285541068, 11.28186362655091, 5.703999834285578, 12.35948216654467, 50.68671012168493, 7.344005520526261, 34.08575148084675, 10.5389719464575, 14.64684517417802)
    25 = @(13.01054333336251, 11.00523134630721, 5.734256494136774, 12.30169429674718, 50.06064657146927, 7.344005520526261, 33.78822663366952, 10.53353860075532, 14.46163248661007)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}
